# Update cryptocurrency price/volume data on the active worksheet.
# Column layout: A=Index, B=Coin, C=Link, D=Price, E=Volume(1h), F=Data, G=Hora

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @{ D = newPrice (optional); E = newVolume (optional) }
$updates = @{
    2  = @{ E = "-0.67%" }
    3  = @{ D = "31.01"; E = "0.91%" }
    4  = @{ D = "4.919"; E = "-0.76%" }
    5  = @{ D = "0.07375"; E = "2.32%" }
    6  = @{ D = "2.213"; E = "24.34%" }
    7  = @{ D = "7.695"; E = "0.30%" }
    8  = @{ D = "3.734"; E = "0.30%" }
    9  = @{ D = "0.9086"; E = "1.43%" }
    10 = @{ D = "0.08743"; E = "13.11%" }
    11 = @{ D = "0.1684"; E = "1.95%" }
    12 = @{ D = "0.08171"; E = "2.17%" }
    13 = @{ D = "0.03115"; E = "2.67%" }
    14 = @{ D = "0.09949"; E = "-0.61%" }
    15 = @{ D = "0.001500"; E = "0.05%" }
    16 = @{ D = "0.005821"; E = "2.58%" }
    17 = @{ D = "3.492"; E = "0.39%" }
    18 = @{ D = "2.065"; E = "-0.87%" }
    19 = @{ E = "0.38%" }
    20 = @{ E = "-1.80%" }
    21 = @{ D = "3.827"; E = "-5.06%" }
    22 = @{ E = "1.10%" }
    23 = @{ D = "0.04553"; E = "0.83%" }
    24 = @{ D = "0.001211"; E = "-0.32%" }
    25 = @{ D = "0.004145"; E = "3.17%" }
    26 = @{ E = "4.16%" }
    27 = @{ D = "0.0003398" }
    39 = @{ D = "0.01581"; E = "-1.28%" }
    40 = @{ D = "0.04467"; E = "1.65%" }
    41 = @{ D = "0.007348"; E = "0.63%" }
    42 = @{ D = "0.009564"; E = "24.12%" }
    43 = @{ E = "1.18%" }
    44 = @{ E = "7.41%" }
    45 = @{ D = "0.008424"; E = "-8.53%" }
    46 = @{ D = "0.00006110"; E = "3.11%" }
    47 = @{ D = "0.00000000751"; E = "0.15%" }
    48 = @{ D = "2.105"; E = "-6.32%" }
    49 = @{ D = "0.002002"; E = "-33.22%" }
    50 = @{ D = "0.00002102"; E = "0.15%" }
    51 = @{ D = "0.0002002"; E = "0.15%" }
}

foreach ($row in $updates.Keys) {
    $cellUpdates = $updates[$row]
    if ($cellUpdates.ContainsKey("D")) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cellUpdates["D"]
    }
    if ($cellUpdates.ContainsKey("E")) {
        $cell = $ws.Range("E$row")
        $cell.NumberFormat = "@"
        $cell.Value = $cellUpdates["E"]
    }
}
